# Update script applied 29-11-2023 02:45
# - Re-shuffles the home/away/odds columns (F:V) for a handful of rows so
#   the data lines up with the matches listed for the correct fixture date
#   (the site re-ordered match listings between scrapes).
# - Appends the newly played match (Fortaleza 1 x 0 Patriotas) as row 156.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 3-way rotation among rows 17, 18, 19 -----------------------------
# final17 = old19, final18 = old17, final19 = old18
$row17 = $ws.Range("F17:V17").Value2
$row18 = $ws.Range("F18:V18").Value2
$row19 = $ws.Range("F19:V19").Value2

$ws.Range("F17:V17").Value2 = $row19
$ws.Range("F18:V18").Value2 = $row17
$ws.Range("F19:V19").Value2 = $row18

# --- simple 2-way swaps -------------------------------------------------
$pairs = @(
    @(44, 45),
    @(67, 68),
    @(77, 78),
    @(92, 94)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    $a = $ws.Range("F$r1`:V$r1").Value2
    $b = $ws.Range("F$r2`:V$r2").Value2
    $ws.Range("F$r1`:V$r1").Value2 = $b
    $ws.Range("F$r2`:V$r2").Value2 = $a
}

# --- append new row 156 --------------------------------------------------
# Copy row 155 so formatting (bold index column, date/time number format)
# carries over, then overwrite with the new match's values. B/C/D
# (pais/torneio/temporada) are identical to row 155 so the copy already
# leaves them correct - re-stamping them with a plain string would make
# "2023" look like a number and pick up a stray style, so leave as-is.
$ws.Range("A155:V155").Copy($ws.Range("A156"))

$ws.Range("A156").Value = 155
$ws.Range("E156").Value = 45259.0625
$ws.Range("F156").Value = "Fortaleza"
$ws.Range("G156").Value = 1
$ws.Range("H156").Value = "Patriotas"
$ws.Range("I156").Value = 0
$ws.Range("J156").Value = 1.69
$ws.Range("K156").Value = "27/11/2023 13:12"
$ws.Range("L156").Value = 1.87
$ws.Range("M156").Value = "29/11/2023 01:28"
$ws.Range("N156").Value = 3.62
$ws.Range("O156").Value = "27/11/2023 13:12"
$ws.Range("P156").Value = 3.32
$ws.Range("Q156").Value = "29/11/2023 01:28"
$ws.Range("R156").Value = 5.28
$ws.Range("S156").Value = "27/11/2023 13:12"
$ws.Range("T156").Value = 4.75
$ws.Range("U156").Value = "29/11/2023 01:28"
$ws.Range("V156").Value = "https://www.betexplorer.com/football/colombia/primera-b/fortaleza-c-e-i-f-patriotas/jeIdXcVu/"
